$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in / clear individual cell values (rows unaffected by later row deletions) ---
$ws.Range("C3").Value = 11.2
$ws.Range("F4").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F13").Value = $null
$ws.Range("F14").Value = $null
$ws.Range("C21").Value = 12.7
$ws.Range("C23").Value = $null

# --- Remove the two rows that no longer appear in the table (RM 232, SC 92) ---
# Row 26 = "RM 232" in the original layout.
$ws.Rows(26).Delete()
# After that deletion, the row that used to be 28 ("SC 92") is now row 27.
$ws.Rows(27).Delete()

# --- SC 193's B column value is now filled in; after both deletions it sits at row 32 ---
$ws.Range("C32").Value = 10.5
